# covid_quest.docx update script
# Applies the wording changes from "confinement" focus to
# "mesures de santé publique" focus, as described by the commit diff.

$d = $word.ActiveDocument
$nbsp = [char]0x00A0

function Replace-Text($findText, $replaceText, [bool]$wildcards) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $wildcards, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING: not found ->" $findText
    }
    return $ok
}

# 1) Title
$f1 = "Évaluation du bien-être durant le confinement en lien avec le COVID-19"
$r1 = "Évaluation du bien-être et lien avec les mesures de santé publique liées au COVID-19"
Replace-Text $f1 $r1 $false

# 2) Intro paragraph - spans a manual line break, so use wildcard match to
#    bridge the two original runs/break and replace with a single new text.
$f2 = "Ce questionnaire étudie les facteurs pouvant avoir un impact sur votre vécu du*confinement en lien avec l'épidémie COVID-19."
$r2 = "Ce questionnaire étudie les facteurs (en particulier : ceux liés aux mesures de santé publique) pouvant avoir un impact sur votre bien-être."
Replace-Text $f2 $r2 $true

# Move the _GoBack bookmark so it lands inside the paragraph we just edited,
# right after "... santé publique" (before ") "), matching the diff.
try {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
} catch {
}
$posRange = $d.Content
$found = $posRange.Find.Execute("ceux liés aux mesures de santé publique", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bmRange = $d.Range($posRange.End, $posRange.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# 3) "Sur une échelle ..." — drop the trailing "depuis le début de la période de confinement"
$f3 = "ces 3 différents domaines depuis le début de la période de confinement" + $nbsp + "?"
$r3 = "ces 3 différents domaines?"
Replace-Text $f3 $r3 $false

# 4) "Coronavirus et confinement" -> "Coronavirus et santé publique"
$f4 = "Coronavirus et confinement"
$r4 = "Coronavirus et santé publique"
Replace-Text $f4 $r4 $false

# 5) "Etes-vous en accord avec la mesure de confinement ? " -> new wording
$f5 = "Etes-vous en accord avec la mesure de confinement" + $nbsp + "? "
$r5 = "Etes-vous en accord avec les mesures actuelles de santé publique liées au COVID-19" + $nbsp + "? "
Replace-Text $f5 $r5 $false

# 6) "... niveau d'information sur les mesures du confinement ? " -> new wording
$f6 = "niveau d’information sur les mesures du confinement" + $nbsp + "? "
$r6 = "niveau d’information sur ces mesures? "
Replace-Text $f6 $r6 $false

# 7) "La période de confinement va-t-elle avoir des répercussions financières sur votre budget ? "
$f7 = "La période de confinement va-t-elle avoir des répercussions financières sur votre budget" + $nbsp + "? "
$r7 = "Les mesures de santé publique actuelles vont-t-elles avoir des répercussions financières sur votre budget" + $nbsp + "? "
Replace-Text $f7 $r7 $false
